$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.593.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.664.28"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4794"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2619"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06154"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07083"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.662.44"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.78"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5923"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.382"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.588.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006742"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.876.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.436"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.661"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.304"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.07"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.46%  "

$ws.Range("E27").Value = "  +1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.79"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.691"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.947"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.662"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07651"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.36%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("E34").Value = "  -4.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.615"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6116"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9500"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.607"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8554"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.000"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.886"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.58%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01502"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.93"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3760"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.710"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1120"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.206"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05262"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.45"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.419"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.04%  "
